# Update master to output generated at 596fc94
# - Bumps the worksheet date heading to the next day
# - Refreshes all 100 two-digit multiplication problems in the table
$d = $word.ActiveDocument

# 1) Date heading (first paragraph, outside the table)
$d.Paragraphs.Item(1).Range.Find.Execute("2023-05-02 Tuesday", $false, $false, $false, $false, $false, $true, 1, $false, "2023-05-03 Wednesday", 2) | Out-Null

# 2) The 20x5 grid of multiplication problems.
#    Cells are addressed directly by (row, column) rather than via
#    Find/Replace so that duplicate-looking intermediate values (e.g. an
#    old problem being replaced by text equal to another still-pending
#    old value) can never cross-match each other.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "60×51="
$t.Cell(1, 2).Range.Text = "57×66="
$t.Cell(1, 3).Range.Text = "68×35="
$t.Cell(1, 4).Range.Text = "46×88="
$t.Cell(1, 5).Range.Text = "65×19="
$t.Cell(2, 1).Range.Text = "54×99="
$t.Cell(2, 2).Range.Text = "10×58="
$t.Cell(2, 3).Range.Text = "57×78="
$t.Cell(2, 4).Range.Text = "42×81="
$t.Cell(2, 5).Range.Text = "14×90="
$t.Cell(3, 1).Range.Text = "37×97="
$t.Cell(3, 2).Range.Text = "59×46="
$t.Cell(3, 3).Range.Text = "84×99="
$t.Cell(3, 4).Range.Text = "87×79="
$t.Cell(3, 5).Range.Text = "45×41="
$t.Cell(4, 1).Range.Text = "13×79="
$t.Cell(4, 2).Range.Text = "49×72="
$t.Cell(4, 3).Range.Text = "52×39="
$t.Cell(4, 4).Range.Text = "57×63="
$t.Cell(4, 5).Range.Text = "91×59="
$t.Cell(5, 1).Range.Text = "16×54="
$t.Cell(5, 2).Range.Text = "93×35="
$t.Cell(5, 3).Range.Text = "50×31="
$t.Cell(5, 4).Range.Text = "37×86="
$t.Cell(5, 5).Range.Text = "97×10="
$t.Cell(6, 1).Range.Text = "81×81="
$t.Cell(6, 2).Range.Text = "61×84="
$t.Cell(6, 3).Range.Text = "27×96="
$t.Cell(6, 4).Range.Text = "59×11="
$t.Cell(6, 5).Range.Text = "53×93="
$t.Cell(7, 1).Range.Text = "64×94="
$t.Cell(7, 2).Range.Text = "32×76="
$t.Cell(7, 3).Range.Text = "53×84="
$t.Cell(7, 4).Range.Text = "82×94="
$t.Cell(7, 5).Range.Text = "84×28="
$t.Cell(8, 1).Range.Text = "92×51="
$t.Cell(8, 2).Range.Text = "90×53="
$t.Cell(8, 3).Range.Text = "68×94="
$t.Cell(8, 4).Range.Text = "49×76="
$t.Cell(8, 5).Range.Text = "66×17="
$t.Cell(9, 1).Range.Text = "27×21="
$t.Cell(9, 2).Range.Text = "91×73="
$t.Cell(9, 3).Range.Text = "25×81="
$t.Cell(9, 4).Range.Text = "40×18="
$t.Cell(9, 5).Range.Text = "19×82="
$t.Cell(10, 1).Range.Text = "21×93="
$t.Cell(10, 2).Range.Text = "50×98="
$t.Cell(10, 3).Range.Text = "94×53="
$t.Cell(10, 4).Range.Text = "89×14="
$t.Cell(10, 5).Range.Text = "98×30="
$t.Cell(11, 1).Range.Text = "35×91="
$t.Cell(11, 2).Range.Text = "17×14="
$t.Cell(11, 3).Range.Text = "63×56="
$t.Cell(11, 4).Range.Text = "70×41="
$t.Cell(11, 5).Range.Text = "11×79="
$t.Cell(12, 1).Range.Text = "57×71="
$t.Cell(12, 2).Range.Text = "25×33="
$t.Cell(12, 3).Range.Text = "30×13="
$t.Cell(12, 4).Range.Text = "55×43="
$t.Cell(12, 5).Range.Text = "38×43="
$t.Cell(13, 1).Range.Text = "45×23="
$t.Cell(13, 2).Range.Text = "92×50="
$t.Cell(13, 3).Range.Text = "55×55="
$t.Cell(13, 4).Range.Text = "61×55="
$t.Cell(13, 5).Range.Text = "99×68="
$t.Cell(14, 1).Range.Text = "46×70="
$t.Cell(14, 2).Range.Text = "57×27="
$t.Cell(14, 3).Range.Text = "84×63="
$t.Cell(14, 4).Range.Text = "100×44="
$t.Cell(14, 5).Range.Text = "25×88="
$t.Cell(15, 1).Range.Text = "22×15="
$t.Cell(15, 2).Range.Text = "41×47="
$t.Cell(15, 3).Range.Text = "11×75="
$t.Cell(15, 4).Range.Text = "61×95="
$t.Cell(15, 5).Range.Text = "42×90="
$t.Cell(16, 1).Range.Text = "66×99="
$t.Cell(16, 2).Range.Text = "13×59="
$t.Cell(16, 3).Range.Text = "74×39="
$t.Cell(16, 4).Range.Text = "93×25="
$t.Cell(16, 5).Range.Text = "58×46="
$t.Cell(17, 1).Range.Text = "78×72="
$t.Cell(17, 2).Range.Text = "66×98="
$t.Cell(17, 3).Range.Text = "80×44="
$t.Cell(17, 4).Range.Text = "13×69="
$t.Cell(17, 5).Range.Text = "92×42="
$t.Cell(18, 1).Range.Text = "77×54="
$t.Cell(18, 2).Range.Text = "48×63="
$t.Cell(18, 3).Range.Text = "71×95="
$t.Cell(18, 4).Range.Text = "78×86="
$t.Cell(18, 5).Range.Text = "58×24="
$t.Cell(19, 1).Range.Text = "47×48="
$t.Cell(19, 2).Range.Text = "26×53="
$t.Cell(19, 3).Range.Text = "32×66="
$t.Cell(19, 4).Range.Text = "95×32="
$t.Cell(19, 5).Range.Text = "68×91="
$t.Cell(20, 1).Range.Text = "76×95="
$t.Cell(20, 2).Range.Text = "25×85="
$t.Cell(20, 3).Range.Text = "13×73="
$t.Cell(20, 4).Range.Text = "24×62="
$t.Cell(20, 5).Range.Text = "35×38="
